$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (data) values first, while rows 3/4 still exist, so the
#     shared-string table gets new entries appended in the same order
#     a person typing the sheet top-to-bottom, left-to-right would produce.
$ws.Range("A2").Value = "Lotte"
$ws.Range("E1").Value = "Trạng thái"
$ws.Range("E2").Value = "x"
$ws.Range("C2").Value = "Xóa-Phế liệu Sáp CRC Xanh lợt"

# B2 is a genuine number.
$ws.Range("B2").Value = 90000631

# --- Formatting copied from existing header/border cells so the new
#     cells reuse the same cell-style indices as their neighbours.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Left-align the new "Mã KH Lotte" number and the text barcode - set the
# alignment BEFORE the number format / value so the barcode cell settles
# on a single combined style instead of an intermediate numFmt-only one.
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("D2").HorizontalAlignment = -4131

# D2 must be stored as TEXT (barcode), so set the number format to Text
# before assigning the value - mirrors typing into a pre-formatted cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "8935001810308"

# --- Drop the old sample rows (3 and 4), now that row 2 holds the new data.
$ws.Rows("3:4").Delete()

# --- Column widths for the new layout.
$ws.Columns("B").ColumnWidth = 15.86
$ws.Columns("C").ColumnWidth = 26.14
$ws.Columns("D").ColumnWidth = 14.86
$ws.Columns("H").ColumnWidth = 9.14

# --- Page setup / selection.
$ws.PageSetup.Orientation = 1
$ws.Range("H13").Select()
